{"js": "// Normalize the placeholder tokens in this identification-letter template so\n// they read as single underscore-joined identifiers (matching the existing\n// {criminal_records} style) instead of containing literal spaces:\n//   {reference }       -> {reference}\n//   {date of birth}    -> {date_of_birth}\n//   {additional info}  -> {additional_info}\n\nconst body = context.document.body;\n\nconst replacements = [\n  { find: \"{reference }\", replace: \"{reference}\" },\n  { find: \"{date of birth}\", replace: \"{date_of_birth}\" },\n  { find: \"{additional info}\", replace: \"{additional_info}\" }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Normalize the placeholder tokens in this identification-letter template so\n# they read as single underscore-joined identifiers (matching the existing\n# {criminal_records} style) instead of containing literal spaces:\n#   {reference }       -> {reference}\n#   {date of birth}    -> {date_of_birth}\n#   {additional info}  -> {additional_info}\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{find = \"{reference }\";      replace = \"{reference}\"},\n    @{find = \"{date of birth}\";   replace = \"{date_of_birth}\"},\n    @{find = \"{additional info}\"; replace = \"{additional_info}\"}\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $p.find\n    $find.Replacement.Text = $p.replace\n    $find.Execute(\n        $p.find,      # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $p.replace,   # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    )\n}\n"}
